# edit.ps1 - applies the homework_2.docx revision:
#   1) Drop the "~" (approx.) prefix from four numeric answers.
#   2) Relocate the "_GoBack" bookmark from the end of the
#      "Yes, by using continuous data improves the model." paragraph
#      to the end of the "TPR: 0.94" paragraph.
#   3) Mark the "Normal Table" table style as a Quick Style (w:qFormat).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Strip the leading "~" from the four approximate values.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Entropy: ~0.93", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Entropy: 0.93", 2) | Out-Null

$d.Content.Find.Execute("Gini coefficient: ~0.47", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Gini coefficient: 0.47", 2) | Out-Null

$d.Content.Find.Execute("Misclassification error: ~0.30", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Misclassification error: 0.30", 2) | Out-Null

$d.Content.Find.Execute("The information gain for the best feature is ~0.21.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "The information gain for the best feature is 0.21.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark to sit right after "TPR: 0.94".
#
#    A bookmark collapsed exactly on the last text offset of a
#    paragraph cannot be created directly in one step here, so a
#    placeholder character is appended after the target text, the
#    bookmark is anchored just before that placeholder, and then the
#    placeholder is removed - leaving the bookmark correctly collapsed
#    right after "TPR: 0.94" and before the paragraph mark.
# ---------------------------------------------------------------------
$target = $d.Content
$found = $target.Find.Execute("TPR: 0.94", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)

if ($found) {
    $anchorPos = $target.End

    # Temporarily append a placeholder right after the found text.
    $target.InsertAfter("X")

    # Build the (now interior, not paragraph-final) collapsed range and
    # drop the bookmark there - this also removes any previous
    # "_GoBack" bookmark elsewhere in the document.
    $bmRange = $d.Range($anchorPos, $anchorPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    # Remove the placeholder character again.
    $placeholder = $d.Range($anchorPos, $anchorPos + 1)
    $placeholder.Delete()
}

# ---------------------------------------------------------------------
# 3) Flag the "Normal Table" table style as a Quick Style.
# ---------------------------------------------------------------------
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true

Write-Output "done"
